$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text is not a plain decimal number (percentages, multi-dot
# "thousand-grouped" prices, subscripted-zero notation, ...): Excel keeps these as
# text automatically, so a plain .Value assignment is enough.
$ws.Range('D2').Value = '27.515.10'
$ws.Range('E2').Value = '  +5.44%  '
$ws.Range('D3').Value = '1.725.38'
$ws.Range('E3').Value = '  +4.72%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +3.16%  '
$ws.Range('E6').Value = '  +2.99%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.66%  '
$ws.Range('E9').Value = '  +4.90%  '
$ws.Range('E10').Value = '  +6.32%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '1.727.20'
$ws.Range('E13').Value = '  +5.54%  '
$ws.Range('D14').Value = '1.962.68'
$ws.Range('E14').Value = '  +4.71%  '
$ws.Range('E15').Value = '  +4.28%  '
$ws.Range('D16').Value = '0.0₅8294'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('E17').Value = '  +4.31%  '
$ws.Range('D18').Value = '27.527.30'
$ws.Range('E18').Value = '  +5.65%  '
$ws.Range('E19').Value = '  +12.77%  '
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +2.94%  '
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('E25').Value = '  -1.03%  '
$ws.Range('E26').Value = '  +14.07%  '
$ws.Range('E27').Value = '  +4.59%  '
$ws.Range('E28').Value = '  +2.12%  '
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('E32').Value = '  +3.37%  '
$ws.Range('E33').Value = '  +3.61%  '
$ws.Range('E34').Value = '  +6.82%  '
$ws.Range('E35').Value = '  +2.79%  '
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('E38').Value = '  +6.33%  '
$ws.Range('E39').Value = '  +4.96%  '
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('D41').Value = '1.048.17'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D45').Value = '1.868.84'
$ws.Range('E45').Value = '  +4.65%  '
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  +2.30%  '
$ws.Range('E47').Value = '  +2.80%  '
$ws.Range('E48').Value = '  +3.77%  '
$ws.Range('E49').Value = '  +3.91%  '
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('E51').Value = '  +2.94%  '

# --- Cells whose new text looks like a plain decimal number (e.g. "225.37").
# A bare .Value assignment would make Excel store these as a Number (and for
# values like "16.50"/"1.750" it would also silently drop the trailing zero).
# Force Text format first so the literal string is preserved exactly, then strip
# the Text-format style back off (paste formats only from an untouched, default
# styled cell) so the cell ends up with its original (default/style-0) formatting
# - only the value itself should change, as in the source edit.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D5').Value = '225.37'
$ws.Range('D6').Value = '0.5357'
$ws.Range('D8').Value = '0.2662'
$ws.Range('D9').Value = '0.06606'
$ws.Range('D10').Value = '21.58'
$ws.Range('D11').Value = '0.07671'
$ws.Range('D12').Value = '4.603'
$ws.Range('D15').Value = '0.5808'
$ws.Range('D17').Value = '67.88'
$ws.Range('D19').Value = '218.85'
$ws.Range('D20').Value = '1.002'
$ws.Range('D21').Value = '4.721'
$ws.Range('D23').Value = '6.037'
$ws.Range('D25').Value = '143.53'
$ws.Range('D26').Value = '1.750'
$ws.Range('D27').Value = '0.1234'
$ws.Range('D28').Value = '7.343'
$ws.Range('D29').Value = '16.50'
$ws.Range('D30').Value = '0.05495'
$ws.Range('D32').Value = '3.550'
$ws.Range('D34').Value = '1.663'
$ws.Range('D36').Value = '0.9580'
$ws.Range('D38').Value = '0.5937'
$ws.Range('D39').Value = '0.01648'
$ws.Range('D40').Value = '5.891'
$ws.Range('D42').Value = '0.8466'
$ws.Range('D43').Value = '1.002'
$ws.Range('D44').Value = '101.30'
$ws.Range('D48').Value = '0.4483'
$ws.Range('D49').Value = '8.176'
$ws.Range('D50').Value = '1.002'
$ws.Range('D51').Value = '0.05248'

$ws.Range("C2").Copy()
$ws.Range('D5').PasteSpecial(-4122)
$ws.Range('D6').PasteSpecial(-4122)
$ws.Range('D8').PasteSpecial(-4122)
$ws.Range('D9').PasteSpecial(-4122)
$ws.Range('D10').PasteSpecial(-4122)
$ws.Range('D11').PasteSpecial(-4122)
$ws.Range('D12').PasteSpecial(-4122)
$ws.Range('D15').PasteSpecial(-4122)
$ws.Range('D17').PasteSpecial(-4122)
$ws.Range('D19').PasteSpecial(-4122)
$ws.Range('D20').PasteSpecial(-4122)
$ws.Range('D21').PasteSpecial(-4122)
$ws.Range('D23').PasteSpecial(-4122)
$ws.Range('D25').PasteSpecial(-4122)
$ws.Range('D26').PasteSpecial(-4122)
$ws.Range('D27').PasteSpecial(-4122)
$ws.Range('D28').PasteSpecial(-4122)
$ws.Range('D29').PasteSpecial(-4122)
$ws.Range('D30').PasteSpecial(-4122)
$ws.Range('D32').PasteSpecial(-4122)
$ws.Range('D34').PasteSpecial(-4122)
$ws.Range('D36').PasteSpecial(-4122)
$ws.Range('D38').PasteSpecial(-4122)
$ws.Range('D39').PasteSpecial(-4122)
$ws.Range('D40').PasteSpecial(-4122)
$ws.Range('D42').PasteSpecial(-4122)
$ws.Range('D43').PasteSpecial(-4122)
$ws.Range('D44').PasteSpecial(-4122)
$ws.Range('D48').PasteSpecial(-4122)
$ws.Range('D49').PasteSpecial(-4122)
$ws.Range('D50').PasteSpecial(-4122)
$ws.Range('D51').PasteSpecial(-4122)
$excel.CutCopyMode = 0

